$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 2.370286
$ws.Cells.Item(2, 8).Value = 7.110858
$ws.Cells.Item(2, 9).Value = 0.3026841782318013
$ws.Cells.Item(2, 10).Value = 0.3026841782318014
$ws.Cells.Item(2, 13).Value = 11.188041
$ws.Cells.Item(2, 14).Value = 33.564123
$ws.Cells.Item(2, 15).Value = 0.1395103797998223
$ws.Cells.Item(2, 16).Value = 0.1395103797998223
$ws.Cells.Item(2, 17).Value = 26.518856949726
$ws.Cells.Item(2, 18).Value = 238.669712547534
$ws.Cells.Item(2, 19).Value = 0.0422275846645157
$ws.Cells.Item(2, 20).Value = 0.04222758466451572

# Row 3
$ws.Cells.Item(3, 7).Value = 2.370286
$ws.Cells.Item(3, 8).Value = 7.110858
$ws.Cells.Item(3, 9).Value = 0.3026841782318013
$ws.Cells.Item(3, 10).Value = 0.3026841782318014
$ws.Cells.Item(3, 15).Value = 0.4168441980730721
$ws.Cells.Item(3, 16).Value = 0.4168441980730722
$ws.Cells.Item(3, 17).Value = 79.23590828785868
$ws.Cells.Item(3, 18).Value = 713.1231745907281
$ws.Cells.Item(3, 19).Value = 0.1261721435444421
$ws.Cells.Item(3, 20).Value = 0.1261721435444421

# Row 4
$ws.Cells.Item(4, 7).Value = 2.370286
$ws.Cells.Item(4, 8).Value = 7.110858
$ws.Cells.Item(4, 9).Value = 0.3026841782318013
$ws.Cells.Item(4, 10).Value = 0.3026841782318014
$ws.Cells.Item(4, 13).Value = 31.78201566666667
$ws.Cells.Item(4, 14).Value = 95.346047
$ws.Cells.Item(4, 15).Value = 0.3963089763847458
$ws.Cells.Item(4, 16).Value = 0.3963089763847459
$ws.Cells.Item(4, 17).Value = 75.33246678648067
$ws.Cells.Item(4, 18).Value = 677.992201078326
$ws.Cells.Item(4, 19).Value = 0.1199564568429032
$ws.Cells.Item(4, 20).Value = 0.1199564568429032

# Row 5
$ws.Cells.Item(5, 7).Value = 2.370286
$ws.Cells.Item(5, 8).Value = 7.110858
$ws.Cells.Item(5, 9).Value = 0.3026841782318013
$ws.Cells.Item(5, 10).Value = 0.3026841782318014
$ws.Cells.Item(5, 13).Value = 3.796148333333333
$ws.Cells.Item(5, 14).Value = 11.388445
$ws.Cells.Item(5, 15).Value = 0.04733644574235969
$ws.Cells.Item(5, 16).Value = 0.04733644574235969
$ws.Cells.Item(5, 17).Value = 8.997957248423333
$ws.Cells.Item(5, 18).Value = 80.98161523581
$ws.Cells.Item(5, 19).Value = 0.01432799317994039
$ws.Cells.Item(5, 20).Value = 0.0143279931799404

# Row 6
$ws.Cells.Item(6, 9).Value = 0.2022126055089961
$ws.Cells.Item(6, 10).Value = 0.2022126055089961
$ws.Cells.Item(6, 13).Value = 11.188041
$ws.Cells.Item(6, 14).Value = 33.564123
$ws.Cells.Item(6, 15).Value = 0.1395103797998223
$ws.Cells.Item(6, 16).Value = 0.1395103797998223
$ws.Cells.Item(6, 17).Value = 17.716311405011
$ws.Cells.Item(6, 18).Value = 159.446802645099
$ws.Cells.Item(6, 19).Value = 0.02821075739487168
$ws.Cells.Item(6, 20).Value = 0.02821075739487169

# Row 7
$ws.Cells.Item(7, 9).Value = 0.2022126055089961
$ws.Cells.Item(7, 10).Value = 0.2022126055089961
$ws.Cells.Item(7, 15).Value = 0.4168441980730721
$ws.Cells.Item(7, 16).Value = 0.4168441980730722
$ws.Cells.Item(7, 19).Value = 0.08429115138366397
$ws.Cells.Item(7, 20).Value = 0.08429115138366398

# Row 8
$ws.Cells.Item(8, 9).Value = 0.2022126055089961
$ws.Cells.Item(8, 10).Value = 0.2022126055089961
$ws.Cells.Item(8, 13).Value = 31.78201566666667
$ws.Cells.Item(8, 14).Value = 95.346047
$ws.Cells.Item(8, 15).Value = 0.3963089763847458
$ws.Cells.Item(8, 16).Value = 0.3963089763847459
$ws.Cells.Item(8, 17).Value = 50.32695953023455
$ws.Cells.Item(8, 18).Value = 452.942635772111
$ws.Cells.Item(8, 19).Value = 0.08013867070136266
$ws.Cells.Item(8, 20).Value = 0.0801386707013627

# Row 9
$ws.Cells.Item(9, 9).Value = 0.2022126055089961
$ws.Cells.Item(9, 10).Value = 0.2022126055089961
$ws.Cells.Item(9, 13).Value = 3.796148333333333
$ws.Cells.Item(9, 14).Value = 11.388445
$ws.Cells.Item(9, 15).Value = 0.04733644574235969
$ws.Cells.Item(9, 16).Value = 0.04733644574235969
$ws.Cells.Item(9, 17).Value = 6.011217335809444
$ws.Cells.Item(9, 18).Value = 54.10095602228499
$ws.Cells.Item(9, 19).Value = 0.009572026029097777
$ws.Cells.Item(9, 20).Value = 0.009572026029097781

# Row 10
$ws.Cells.Item(10, 7).Value = 2.286703333333333
$ws.Cells.Item(10, 8).Value = 6.860109999999999
$ws.Cells.Item(10, 9).Value = 0.2920107190904054
$ws.Cells.Item(10, 10).Value = 0.2920107190904054
$ws.Cells.Item(10, 13).Value = 11.188041
$ws.Cells.Item(10, 14).Value = 33.564123
$ws.Cells.Item(10, 15).Value = 0.1395103797998223
$ws.Cells.Item(10, 16).Value = 0.1395103797998223
$ws.Cells.Item(10, 17).Value = 25.58373064817
$ws.Cells.Item(10, 18).Value = 230.25357583353
$ws.Cells.Item(10, 19).Value = 0.04073852632592168
$ws.Cells.Item(10, 20).Value = 0.04073852632592168

# Row 11
$ws.Cells.Item(11, 7).Value = 2.286703333333333
$ws.Cells.Item(11, 8).Value = 6.860109999999999
$ws.Cells.Item(11, 9).Value = 0.2920107190904054
$ws.Cells.Item(11, 10).Value = 0.2920107190904054
$ws.Cells.Item(11, 15).Value = 0.4168441980730721
$ws.Cells.Item(11, 16).Value = 0.4168441980730722
$ws.Cells.Item(11, 17).Value = 76.44183680852889
$ws.Cells.Item(11, 18).Value = 687.9765312767599
$ws.Cells.Item(11, 19).Value = 0.1217229740279812
$ws.Cells.Item(11, 20).Value = 0.1217229740279812

# Row 12
$ws.Cells.Item(12, 7).Value = 2.286703333333333
$ws.Cells.Item(12, 8).Value = 6.860109999999999
$ws.Cells.Item(12, 9).Value = 0.2920107190904054
$ws.Cells.Item(12, 10).Value = 0.2920107190904054
$ws.Cells.Item(12, 13).Value = 31.78201566666667
$ws.Cells.Item(12, 14).Value = 95.346047
$ws.Cells.Item(12, 15).Value = 0.3963089763847458
$ws.Cells.Item(12, 16).Value = 0.3963089763847459
$ws.Cells.Item(12, 17).Value = 72.67604116501889
$ws.Cells.Item(12, 18).Value = 654.0843704851699
$ws.Cells.Item(12, 19).Value = 0.1157264691760921
$ws.Cells.Item(12, 20).Value = 0.1157264691760922

# Row 13
$ws.Cells.Item(13, 7).Value = 2.286703333333333
$ws.Cells.Item(13, 8).Value = 6.860109999999999
$ws.Cells.Item(13, 9).Value = 0.2920107190904054
$ws.Cells.Item(13, 10).Value = 0.2920107190904054
$ws.Cells.Item(13, 13).Value = 3.796148333333333
$ws.Cells.Item(13, 14).Value = 11.388445
$ws.Cells.Item(13, 15).Value = 0.04733644574235969
$ws.Cells.Item(13, 16).Value = 0.04733644574235969
$ws.Cells.Item(13, 17).Value = 8.68066504766111
$ws.Cells.Item(13, 18).Value = 78.12598542894997
$ws.Cells.Item(13, 19).Value = 0.01382274956041041
$ws.Cells.Item(13, 20).Value = 0.01382274956041041

# Row 14
$ws.Cells.Item(14, 7).Value = 1.590394666666667
$ws.Cells.Item(14, 8).Value = 4.771184
$ws.Cells.Item(14, 9).Value = 0.2030924971687972
$ws.Cells.Item(14, 10).Value = 0.2030924971687972
$ws.Cells.Item(14, 13).Value = 11.188041
$ws.Cells.Item(14, 14).Value = 33.564123
$ws.Cells.Item(14, 15).Value = 0.1395103797998223
$ws.Cells.Item(14, 16).Value = 0.1395103797998223
$ws.Cells.Item(14, 17).Value = 17.793400736848
$ws.Cells.Item(14, 18).Value = 160.140606631632
$ws.Cells.Item(14, 19).Value = 0.02833351141451322
$ws.Cells.Item(14, 20).Value = 0.02833351141451323

# Row 15
$ws.Cells.Item(15, 7).Value = 1.590394666666667
$ws.Cells.Item(15, 8).Value = 4.771184
$ws.Cells.Item(15, 9).Value = 0.2030924971687972
$ws.Cells.Item(15, 10).Value = 0.2030924971687972
$ws.Cells.Item(15, 15).Value = 0.4168441980730721
$ws.Cells.Item(15, 16).Value = 0.4168441980730722
$ws.Cells.Item(15, 17).Value = 53.16504672832712
$ws.Cells.Item(15, 18).Value = 478.485420554944
$ws.Cells.Item(15, 19).Value = 0.08465792911698493
$ws.Cells.Item(15, 20).Value = 0.08465792911698494

# Row 16
$ws.Cells.Item(16, 7).Value = 1.590394666666667
$ws.Cells.Item(16, 8).Value = 4.771184
$ws.Cells.Item(16, 9).Value = 0.2030924971687972
$ws.Cells.Item(16, 10).Value = 0.2030924971687972
$ws.Cells.Item(16, 13).Value = 31.78201566666667
$ws.Cells.Item(16, 14).Value = 95.346047
$ws.Cells.Item(16, 15).Value = 0.3963089763847458
$ws.Cells.Item(16, 16).Value = 0.3963089763847459
$ws.Cells.Item(16, 17).Value = 50.54594821218311
$ws.Cells.Item(16, 18).Value = 454.913533909648
$ws.Cells.Item(16, 19).Value = 0.0804873796643879
$ws.Cells.Item(16, 20).Value = 0.08048737966438792

# Row 17
$ws.Cells.Item(17, 7).Value = 1.590394666666667
$ws.Cells.Item(17, 8).Value = 4.771184
$ws.Cells.Item(17, 9).Value = 0.2030924971687972
$ws.Cells.Item(17, 10).Value = 0.2030924971687972
$ws.Cells.Item(17, 13).Value = 3.796148333333333
$ws.Cells.Item(17, 14).Value = 11.388445
$ws.Cells.Item(17, 15).Value = 0.04733644574235969
$ws.Cells.Item(17, 16).Value = 0.04733644574235969
$ws.Cells.Item(17, 17).Value = 6.037374063208889
$ws.Cells.Item(17, 18).Value = 54.33636656888
$ws.Cells.Item(17, 19).Value = 0.009613676972911106
$ws.Cells.Item(17, 20).Value = 0.009613676972911107

